# myapp DAO, observer pattern
#
# Append the newly-seeded rows (user7/user8, board posts "ee"/"ff", and
# project "p6") to the existing "users", "boards" and "projects" sheets.
#
# Purely-numeric-looking fields (no/password/view_count/start_date/...)
# are written through a small helper that forces them to stay TEXT
# (matching every other cell in these sheets, which are all shared
# strings) instead of letting Excel auto-coerce them to Number, while
# restoring the default (General) number format afterwards so no new
# cell style gets introduced.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- users: add user7 (row 7) and user8 (row 8) ---
$usersSheet = $wb.Worksheets.Item("users")

Set-TextValue $usersSheet.Range("A7") "7"
$usersSheet.Range("B7").Value = "user7"
$usersSheet.Range("C7").Value = "user7@test.co"
Set-TextValue $usersSheet.Range("D7") "777"
$usersSheet.Range("E7").Value = "010-1111-7777"

Set-TextValue $usersSheet.Range("A8") "8"
$usersSheet.Range("B8").Value = "user8"
$usersSheet.Range("C8").Value = "user8@test.com"
Set-TextValue $usersSheet.Range("D8") "8888"
$usersSheet.Range("E8").Value = "010-1111-8888"

# --- boards: add board "ee" (row 5) and "ff" (row 6) ---
$boardsSheet = $wb.Worksheets.Item("boards")

Set-TextValue $boardsSheet.Range("A5") "5"
$boardsSheet.Range("B5").Value = "ee"
$boardsSheet.Range("C5").Value = "ee"
$boardsSheet.Range("D5").Value = "2024-07-26 11:21:36"
Set-TextValue $boardsSheet.Range("E5") "0"

Set-TextValue $boardsSheet.Range("A6") "6"
$boardsSheet.Range("B6").Value = "ff"
$boardsSheet.Range("C6").Value = "ff"
$boardsSheet.Range("D6").Value = "2024-07-26 16:46:35"
Set-TextValue $boardsSheet.Range("E6") "0"

# --- projects: add project "p6" (row 5) ---
$projectsSheet = $wb.Worksheets.Item("projects")

Set-TextValue $projectsSheet.Range("A5") "5"
$projectsSheet.Range("B5").Value = "p6"
$projectsSheet.Range("C5").Value = "xx"
Set-TextValue $projectsSheet.Range("D5") "66"
Set-TextValue $projectsSheet.Range("E5") "77"
Set-TextValue $projectsSheet.Range("F5") "5"
